$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("USUARIOS")

# Row 2 values change from (1,5,1,1) to (2,2,2,2). The original cells store
# these as text, not numbers, so force a text number-format while writing
# them, then clear the format again so no stray style sticks around.
$ws.Range("A2:D2").NumberFormat = "@"
$ws.Range("A2").Value = "2"
$ws.Range("B2").Value = "2"
$ws.Range("C2").Value = "2"
$ws.Range("D2").Value = "2"
$ws.Range("A2:D2").ClearFormats()

# Rows 4-9 are removed entirely, shrinking the used range down to A1:D3.
$ws.Range("A4:D9").EntireRow.Delete()
